# Auto-update draw results: append the 2025-10-19 "Pick 4" draw as a new
# row (row 33) at the bottom of the Results sheet, mirroring the existing
# table layout: Date | Game | Phase | Result | InsertedAt.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 33

# Columns A (Date, e.g. "2025-10-19") and C (Phase, e.g. "251019") look
# like numbers/dates, so format them as Text first to keep them as plain
# string values instead of being auto-converted to a date serial / number
# (matching every other row already in the sheet).
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 3).NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = "2025-10-19"
$ws.Cells.Item($newRow, 2).Value = "Pick 4"
$ws.Cells.Item($newRow, 3).Value = "251019"
$ws.Cells.Item($newRow, 4).Value = "4-6-0-9"
$ws.Cells.Item($newRow, 5).Value = "2025-10-19T21:36:17.364+04:00"
